$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "instruction"
$ws.Range("D3").Value = "expr"
$ws.Range("C3").Value = "id"
$ws.Range("B2").Value = "name"
$ws.Range("B3").Value = "name"

$ws.Range("D4").Select()
$excel.ActiveWindow.Zoom = 142
